$d = $word.ActiveDocument

# --- Main body (document.xml): "A TERE," -> "A QWER," ---
$d.Content.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- Header (header1.xml): replace each TRE/TERE/Tre/tre occurrence, in document order,
#     with its corresponding new value from the diff. ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)

$headerTargets = @(
  @("TRE", "QWER"),
  @("TERE", "QWER"),
  @("Tre", "Qwer"),
  @("Tre", "Qwer"),
  @("Tre", "Qewr"),
  @("Tre", "Qewr"),
  @("Tre", "Qwer"),
  @("tre", "qwer"),
  @("tre", "qwer"),
  @("tre", "qwer")
)

foreach ($pair in $headerTargets) {
    $old = $pair[0]
    $new = $pair[1]
    $r = $hdr.Range
    $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}
